$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.535.43'
$ws.Range('E2').Value = '  +2.50%  '
$ws.Range('D3').Value = '2.470.26'
$ws.Range('E3').Value = '  +2.10%  '
$ws.Range('E4').Value = '  +0.24%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '576.21'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +2.39%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '148.22'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +3.59%  '
$ws.Range('E7').Value = '  -0.08%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.541'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +1.82%  '
$ws.Range('E9').Value = '  +4.32%  '
$ws.Range('E10').Value = '  +0.75%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '5.34'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +2.73%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.363'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +3.82%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '27.26'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +3.85%  '
$ws.Range('E14').Value = '  +6.58%  '
$ws.Range('D15').Value = '2.911.18'
$ws.Range('E15').Value = '  +2.34%  '
$ws.Range('D16').Value = '63.463.28'
$ws.Range('E16').Value = '  +2.53%  '
$ws.Range('D17').Value = '2.476.18'
$ws.Range('E17').Value = '  +2.26%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '11.55'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +1.89%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '7.30'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +6.99%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '4.25'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +2.61%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '328.67'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +1.51%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '0.997'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('E23').Value = '  +11.21%  '
$ws.Range('E24').Value = '  +0.83%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '626.99'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +13.34%  '
$ws.Range('E26').Value = '  +13.09%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '8.75'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('E28').Value = '  +2.27%  '
$ws.Range('E29').Value = '  +9.54%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '8.45'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +2.70%  '
$ws.Range('B31').Value = 'Binance-PegBSC-USD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -0.19%  '
$ws.Range('E32').Value = '  -1.47%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '1.93'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +2.77%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '5.19'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +9.70%  '
$ws.Range('E35').Value = '  +3.47%  '
$ws.Range('E36').Value = '  -0.18%  '
$ws.Range('E37').Value = '  +1.94%  '
$ws.Range('E38').Value = '  +1.61%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '18.98'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +1.59%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '1.86'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +2.51%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '147.25'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -3.40%  '
$ws.Range('E42').Value = '  +19.87%  '
$ws.Range('E43').Value = '  +0.59%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '150.36'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +2.04%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '3.77'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +3.58%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.0550'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +4.21%  '
$ws.Range('E47').Value = '  +6.59%  '
$ws.Range('E48').Value = '  +2.41%  '
$ws.Range('E49').Value = '  +5.60%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.0928'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +0.93%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.747'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +4.73%  '
